$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.436.80"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.852.14"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.23"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.48%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4745"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.35%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2752"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.07%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06324"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "17.65"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +9.83%  "
$ws.Range("D11").Value = "1.877.11"
$ws.Range("E11").Value = "  +2.39%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07455"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.952"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "84.75"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.12%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6255"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "30.404.63"
$ws.Range("E16").Value = "  +1.03%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "246.24"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +7.71%  "
$ws.Range("E18").Value = "  +0.03%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.68"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007330"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.908"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.903"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.57%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "164.97"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.096"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +1.53%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.873"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.98%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.1031"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.87%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.347"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.55%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.041"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.825"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.13%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.04852"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +0.64%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.6989"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  +4.98%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.680"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.84%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.8784"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.44%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.997"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.41%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "106.75"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.52%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4060"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("E43").Value = "  +0.51%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "7.173"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.90%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "63.40"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.95%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.1198"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "33.93"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.62%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.543"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  -0.29%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.351"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3692"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.99%  "
